# Refresh the "Generated" timestamps in the Ecom Operations Tracking System
# workbook. The report generation timestamp moved forward by ~4 minutes
# (from 12:10:45 to 12:14:13), and every "last updated"/timestamp cell that
# was derived from that same generation moment shifts from :10 to :14
# accordingly (minutes component only changes, seconds/hours stay aligned
# with each cell's own original value).

$wb = $excel.ActiveWorkbook

# --- Dashboard sheet -------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("A2").Value = "Generated: 2025-10-06 12:14:13"
$dash.Range("E6").Value  = "2025-10-06 12:14"
$dash.Range("E7").Value  = "2025-10-06 12:14"
$dash.Range("E8").Value  = "2025-10-06 12:14"
$dash.Range("E9").Value  = "2025-10-06 12:14"
$dash.Range("E10").Value = "2025-10-06 12:14"
$dash.Range("E11").Value = "2025-10-06 12:14"
$dash.Range("E12").Value = "2025-10-06 12:14"

# --- Bash Queries Response sheet -------------------------------------
$bash = $wb.Worksheets.Item("Bash Queries Response")
$bash.Range("B3").Value = "2025-10-06 12:14"
$bash.Range("B4").Value = "2025-10-06 12:14"
$bash.Range("B5").Value = "2025-10-06 12:14"

# --- System Errors sheet ----------------------------------------------
$errs = $wb.Worksheets.Item("System Errors")
$errs.Range("B3").Value = "2025-10-06 12:14"
$errs.Range("B4").Value = "2025-10-06 12:14"
$errs.Range("B5").Value = "2025-10-06 12:14"

# --- Stock Replenishment sheet -----------------------------------------
$stock = $wb.Worksheets.Item("Stock Replenishment")
$stock.Range("G3").Value = "2025-10-06 09:14"
$stock.Range("H3").Value = "2025-10-06 11:14"
$stock.Range("G4").Value = "2025-10-06 07:14"
$stock.Range("H4").Value = "2025-10-06 12:14"
$stock.Range("G5").Value = "2025-10-06 10:14"
